$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Objetivos:" body text becomes the docente text
$ws.Range("B10").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C10").Value = "1285870 - Marcos Villela Barcza"

# Row 13: gains label "Programa resumido:" in A13 (new cell - copy format from A12
# first so it carries the correct column-A style), and B13/C13 become "Semestral"
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: becomes "Short syllabus:" only label, B14/C14 cleared entirely
# (use Clear(), not Delete(), so nothing below shifts up)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# Row 15: becomes "Programa:" label; B15/C15 (new cells) get "01/01/2016" - copy
# from B8/C8 (already holding that exact text as a shared string) so no date
# auto-conversion or new number-format style gets introduced
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows.Item(15).RowHeight = 120

# Row 16: becomes "Syllabus:" label only, B16/C16 cleared entirely
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Rows.Item(16).RowHeight = 120

# Row 17: becomes "Avaliação:" label only, default row height (no custom height)
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# Row 18: becomes "Método:"; B18/C18 (new cells) get the docente text - copy
# format from B19/C19 first so the correct body style is used
$ws.Range("A18").Value = "Método:"
$ws.Range("B19").Copy($ws.Range("B18"))
$ws.Range("C19").Copy($ws.Range("C18"))
$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: becomes "Critério:" label (its body text in B19/C19 is unchanged)
$ws.Range("A19").Value = "Critério:"

# Row 20: becomes "Norma de recuperação:" label (its body text unchanged)
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: becomes "Bibliografia:" label (its body text unchanged), row grows to 120
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22: becomes "Requisitos:" label only, B22/C22 cleared, default row height
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# Row 23: A23 removed (single-column delete only affects column A, nothing below
# it in that column, so this is safe); B23/C23 (new cells) get the requisito
# text - copy format from B19/C19 (correct body style) first
$ws.Range("A23").Delete()
$ws.Range("B19").Copy($ws.Range("B23"))
$ws.Range("C19").Copy($ws.Range("C23"))
$ws.Range("B23").Value = "LOQ4047 -  Trabalho de Conclusão de Curso I  (Requisito)`n"
$ws.Range("C23").Value = "LOQ4047 -  Trabalho de Conclusão de Curso I  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# Former row 24 (now unused beyond the new extent) is removed entirely - it is
# the last row, so this full-row delete cannot shift any other data
$ws.Rows.Item(24).Delete()
